$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("15").Insert()
